$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new, more realistic demo values
$ws.Range("A2").Value = "ser_pub_loc___variable_3"
$ws.Range("B2").Value = "modality_1"
$ws.Range("A3").Value = "ser_pub_loc___variable_5"
$ws.Range("B3").Value = "modality_1"
$ws.Range("A4").Value = "accident_route___variable_2"
$ws.Range("B4").Value = "a_or_b"

# Add new demo rows 5-8
$ws.Range("A5").Value = "ser_pub_loc___canton"
$ws.Range("B5").Value = "canton_sigle"
$ws.Range("A6").Value = "ser_pub_loc___langue"
$ws.Range("B6").Value = "langue_sigle"
$ws.Range("A7").Value = "ser_pub_loc___nouveau"
$ws.Range("B7").Value = "oui_non"
$ws.Range("A8").Value = "ser_pub_loc___nouveau"
$ws.Range("B8").Value = "vide"

# Resize the table to cover the newly added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B8"))

# Update the active selection to reflect where the user ended up after editing
$ws.Range("B9").Select()
